$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 117.75
$ws.Range("I5").Value = 117.75
$ws.Range("K5").Value = 117.75
$ws.Range("M5").Value = -2.75

$ws.Range("H40").Value = 5998.3335
$ws.Range("J40").Value = 7500
$ws.Range("L40").Value = 7500
$ws.Range("N40").Value = -7850

$ws.Range("H62").Value = 7000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 7000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -8126
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -22632
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 2750
$ws.Range("I80").Value = 2750
$ws.Range("K80").Value = 8250
$ws.Range("M80").Value = -7252

$ws.Range("H83").Value = 2750
$ws.Range("I83").Value = 2750
$ws.Range("K83").Value = 24750
$ws.Range("M83").Value = -19758

$ws.Range("H100").Value = 1495.1111
$ws.Range("I100").Value = 1434.6666
$ws.Range("J100").Value = 1616
$ws.Range("K100").Value = 1434.6666
$ws.Range("L100").Value = 1616
$ws.Range("M100").Value = -893.6666
$ws.Range("N100").Value = -2698

$ws.Range("H135").Value = 894.6667
$ws.Range("I135").Value = 346.4
$ws.Range("K135").Value = 3117.6
$ws.Range("M135").Value = -582.5999999999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2285.25
$ws.Range("I61").Value = 1175.5834
$ws.Range("K61").Value = 1175.5834
$ws.Range("M61").Value = -963.5834

$ws.Range("H74").Value = 858.75
$ws.Range("I74").Value = 858.75
$ws.Range("K74").Value = 858.75
$ws.Range("M74").Value = 15.25

$ws.Range("H77").Value = 858.75
$ws.Range("I77").Value = 858.75
$ws.Range("K77").Value = 4293.75
$ws.Range("M77").Value = 74.25

$ws.Range("H97").Value = 5353.2856
$ws.Range("I97").Value = 4745.3335
$ws.Range("K97").Value = 4745.3335
$ws.Range("M97").Value = -4249.3335

$ws.Range("H132").Value = 2731.0476
$ws.Range("I132").Value = 1335.3077
$ws.Range("J132").Value = 4999.125
$ws.Range("K132").Value = 4005.9231
$ws.Range("L132").Value = 14997.375
$ws.Range("M132").Value = -1475.9231
$ws.Range("N132").Value = -20057.375

$ws.Range("H136").Value = 2285.25
$ws.Range("I136").Value = 1175.5834
$ws.Range("K136").Value = 3526.7502
$ws.Range("M136").Value = -976.7502


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2021.1666
$ws.Range("I86").Value = 1906
$ws.Range("K86").Value = 1906
$ws.Range("M86").Value = -783

$ws.Range("H89").Value = 2021.1666
$ws.Range("I89").Value = 1906
$ws.Range("K89").Value = 9530
$ws.Range("M89").Value = -3914


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3390.5757
$ws.Range("I31").Value = 2578.75
$ws.Range("J31").Value = 5555.4443
$ws.Range("K31").Value = 2578.75
$ws.Range("L31").Value = 5555.4443
$ws.Range("M31").Value = -2283.75
$ws.Range("N31").Value = -6145.4443

$ws.Range("H34").Value = 3390.5757
$ws.Range("I34").Value = 2578.75
$ws.Range("J34").Value = 5555.4443
$ws.Range("K34").Value = 2578.75
$ws.Range("L34").Value = 5555.4443
$ws.Range("M34").Value = -2376.75
$ws.Range("N34").Value = -5959.4443


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 10000
$ws.Range("J138").Value = 10000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 29545
$ws.Range("J54").Value = 29545
$ws.Range("L54").Value = 29545
$ws.Range("N54").Value = -30325

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H97").Value = 2216.6667
$ws.Range("I97").Value = 1750
$ws.Range("J97").Value = 3150
$ws.Range("K97").Value = 1750
$ws.Range("L97").Value = 3150
$ws.Range("M97").Value = -1254
$ws.Range("N97").Value = -4142

$ws.Range("H99").Value = 39950
$ws.Range("I99").Value = 39950
$ws.Range("K99").Value = 39950
$ws.Range("M99").Value = -37704

$ws.Range("H132").Value = 2118.611
$ws.Range("I132").Value = 1731.2142
$ws.Range("J132").Value = 3474.5
$ws.Range("K132").Value = 5193.642599999999
$ws.Range("L132").Value = 10423.5
$ws.Range("M132").Value = -2663.642599999999
$ws.Range("N132").Value = -15483.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3080.4
$ws.Range("I40").Value = 3080.4
$ws.Range("K40").Value = 3080.4
$ws.Range("M40").Value = -2944.4

$ws.Range("H132").Value = 3286.5789
$ws.Range("I132").Value = 3058.7273
$ws.Range("J132").Value = 3599.875
$ws.Range("K132").Value = 9176.1819
$ws.Range("L132").Value = 10799.625
$ws.Range("M132").Value = -6646.1819
$ws.Range("N132").Value = -15859.625


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5888
$ws.Range("I62").Value = 5888
$ws.Range("K62").Value = 5888
$ws.Range("M62").Value = -5264

$ws.Range("H65").Value = 5888
$ws.Range("I65").Value = 5888
$ws.Range("K65").Value = 29440
$ws.Range("M65").Value = -26320

$ws.Range("H126").Value = 1934
$ws.Range("I126").Value = 901
$ws.Range("K126").Value = 2703
$ws.Range("M126").Value = -233

$ws.Range("H136").Value = 2297.5
$ws.Range("I136").Value = 2300
$ws.Range("J136").Value = 2295
$ws.Range("K136").Value = 6900
$ws.Range("L136").Value = 6885
$ws.Range("M136").Value = -4350
$ws.Range("N136").Value = -11985

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

